$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Remove existing hyperlinks on F2:F13 before the rows shift down,
#    otherwise the hyperlink <-> cell association breaks on insert.
$ws.Range("F2:F13").Hyperlinks.Delete()

# 2) Insert 7 fresh rows above the current row 2, pushing the existing
#    13 data rows down to rows 9-20.
$ws.Rows("2:8").Insert()

# 3) Populate the 7 newly inserted rows with the new listings.
$ws.Range("B2").Value = '医療機関向けAIアプリとLINEの連携開発を支援してくださるAIエンジニア募集(AI/バックエンド)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = '🔥AI,Ai ◆開発 ◇アプリ'

$ws.Range("B3").Value = '大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("G3").Value = 378
$ws.Range("H3").Value = '🔥AI,Ai ◆効率化'

$ws.Range("B4").Value = 'Azureでの社内文書検索RAG開発の精度改善を伴走支援してくださるAIエンジニア募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

$ws.Range("B5").Value = 'Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("G5").Value = 310
$ws.Range("H5").Value = '🔥AI,Ai'

$ws.Range("B6").Value = 'Webシステム チャット機能へのChatwork連携API新規開発・組み込み'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("G6").Value = 265
$ws.Range("H6").Value = '🔥API ◆開発'

$ws.Range("B7").Value = '【急募】PDF見積書をExcel注文書に変換するシステム開発'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("G7").Value = 118
$ws.Range("H7").Value = '◆開発,システム開発'

$ws.Range("B8").Value = '【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("G8").Value = 115
$ws.Range("H8").Value = '◆開発 ◇アプリ'

# 4) Refresh the "fetched at" timestamp for every row (new + shifted).
$timestamp = '2025-10-27 18:26:06'
for ($r = 2; $r -le 20; $r++) {
    $ws.Range("A$r").Value = $timestamp
}

# 5) Re-create the hyperlinks for every URL cell F2:F20 (new + shifted),
#    applying the built-in Hyperlink style so formatting matches.
$urls = @{
    2 = 'https://www.lancers.jp/work/detail/5416301'
    3 = 'https://www.lancers.jp/work/detail/5416307'
    4 = 'https://www.lancers.jp/work/detail/5416305'
    5 = 'https://www.lancers.jp/work/detail/5416328'
    6 = 'https://www.lancers.jp/work/detail/5421556'
    7 = 'https://www.lancers.jp/work/detail/5421561'
    8 = 'https://www.lancers.jp/work/detail/5421687'
    9 = 'https://www.lancers.jp/work/detail/5421083'
    10 = 'https://www.lancers.jp/work/detail/5421443'
    11 = 'https://www.lancers.jp/work/detail/5421265'
    12 = 'https://www.lancers.jp/work/detail/5420868'
    13 = 'https://www.lancers.jp/work/detail/5421445'
    14 = 'https://www.lancers.jp/work/detail/5421105'
    15 = 'https://www.lancers.jp/work/detail/5420971'
    16 = 'https://www.lancers.jp/work/detail/5421230'
    17 = 'https://www.lancers.jp/work/detail/5341051'
    18 = 'https://www.lancers.jp/work/detail/5421564'
    19 = 'https://www.lancers.jp/work/detail/5421418'
    20 = 'https://www.lancers.jp/work/detail/5421177'
}
foreach ($r in 2..20) {
    $cell = $ws.Range("F$r")
    $ws.Hyperlinks.Add($cell, $urls[$r]) | Out-Null
    $cell.Style = "Hyperlink"
}

# 6) Widen column B (46 -> 52 chars). COM's ColumnWidth setter pads by
#    5/6 of a character internally, so back it off to land on 52 exactly.
$ws.Columns("B:B").ColumnWidth = 52 - (5/6)
